$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $ok = $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $ok) {
        Write-Output "WARNING: replacement not found for: $find"
    }
}

# 1. Title heading and the later bold duplicate (both occurrences get replaced by Replace:=2 / wdReplaceAll)
Replace-Text "Play Cleopatra Diamond Spins Free | Review of IGT's Slot Game" "Play Cleopatra Diamond Spins for Free"

# 2. "What we like" bullet list
Replace-Text "Easy to learn and play for all types of players" "Easy to play"
Replace-Text "High-quality graphics and well-designed symbols" "Multiple device compatibility"
Replace-Text "Access through multiple devices including smartphones and tablets" "High-quality graphics"
Replace-Text "Immersive experience created by sound and graphics" "Immersive sound and music"

# 3. "What we don't like" bullet list
Replace-Text "Does not offer much innovation compared to other Egyptian themed slot games" "Not unique in its theme"
Replace-Text "Limited special functions and bonuses" "Limited special features"

# 4. Meta description italic text at the end
Replace-Text "Learn about IGT's Cleopatra Diamond Spins slot game in this review, and play for free. Features, symbols, gameplay, and sound and graphics are discussed." "Discover the pros and cons of Cleopatra Diamond Spins and play for free today."
